$p = $ppt.ActivePresentation

# Slide 15: "Some Open-Source Software Resources" -> item 2 (Python Scientific Stack -> pydata)
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 4: "2. Python Scientific Stack" -> "2. Python " + "Open Data Science Stack"
$para4 = $tr.Paragraphs(4)
$r = $para4.Runs(1)
$r.Text = "2. Python "
$r2 = $para4.InsertAfter("Open Data Science Stack")

# Paragraph 5: "https://www.scipy.org/index.html" -> "http://" + "pydata.org/downloads.html"
# (hyperlink target relationship rId3 is left untouched, matching original scipy.org URL)
$para5 = $tr.Paragraphs(5)
$r1b = $para5.Runs(1)
$r1b.Text = "http://pydata.org/downloads.html"
$c2 = $r1b.Characters(8, 25)
$acts = $c2.ActionSettings
$a = $acts.Item(1)
$h = $a.Hyperlink
$h.Address = "https://www.scipy.org/index.html"

Write-Host "done"
